$d = $word.ActiveDocument

# Step 1: merge the "<<cs_{" / "writtenByJudge" / "}>><<" / "hearingLocation." runs
# (and drop their surrounding spell-check proofErr markers) into a single run.
$d.Content.Find.Execute("<<cs_{writtenByJudge}>><<hearingLocation.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<<cs_{writtenByJudge}>><<hearingLocation.", 2)

# Step 2: swap the court "venue" placeholder for the "external_short" placeholder.
$d.Content.Find.Execute("venue", $true, $false, $false, $false, $false,
                         $true, 1, $false, "external_short", 2)

# Step 3: merge the trailing "_name" / ">><<else>> Online Civil Claims<<es_>>" runs into a
# single run, and drop the stray space that used to precede "Online Civil Claims".
$d.Content.Find.Execute("_name>><<else>> Online Civil Claims<<es_>>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "_name>><<else>>Online Civil Claims<<es_>>", 2)
